# add Session 16 and 17
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert 7 blank rows before the trailing "AOP/Security/Microservice/Cloud"
# block (currently rows 98-101) so that block lands on rows 105-108, and the
# new Session 16 / Session 17 / Session 18 content can be written into the
# freshly opened rows 92-104.
$ws.Rows("92:98").Insert()

# Give the new section header bars (A92, A95, A102) the same look as the
# other "Session N" header bars on the sheet (merged A:E, bold/fill style).
$ws.Range("A84:E84").Copy()
$ws.Range("A92:E92").PasteSpecial(-4122)
$ws.Range("A95:E95").PasteSpecial(-4122)
$ws.Range("A102:E102").PasteSpecial(-4122)
$ws.Range("A92:E92").Merge()
$ws.Range("A95:E95").Merge()
$ws.Range("A102:E102").Merge()

# ---- Session 16 ----
$ws.Range("A92").Value = "Session 16"
$ws.Range("D93").Value = "Support XML and JSON for Request and Response with accept header"
$ws.Range("D94").Value = "Update User"

# ---- Session 17 ----
$ws.Range("A95").Value = "Session 17"
$ws.Range("D96").Value = "Make Custom Exception and CustomErrors"
$ws.Range("D97").Value = "Write ControllerAdvice for handle errors"
$ws.Range("D98").Value = "Results Pagination for getAll using PagingAndSortingRepository"
$ws.Range("D99").Value = "create contextpath and Build jar file and "
$ws.Range("D100").Value = "Create war file - SpringBootServletInitializer, tomcat as provider"
$ws.Range("D101").Value = "make tomcat server"

# New rows under the (now shifted) AOP/Security/Microservice/Cloud block
$ws.Range("B110").Value = "SOAP"
$ws.Range("B110").Interior.ThemeColor = 6
$ws.Range("B111").Value = "Rest Client"

# ---- Session 18 ----
$ws.Range("D103").Value = "RestClient"
$ws.Range("D104").Value = "SOAP Server And Client"
$ws.Range("A102").Value = "Session 18"

# Restore the view roughly where the author left it
[void]$ws.Range("A103").Select()
$excel.ActiveWindow.ScrollRow = 84
